$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("T17").Value = 1
$ws.Range("T23").Value = 1
$ws.Range("T24").Value = 1
$ws.Range("T25").Value = 1
$ws.Range("T26").Value = 1
$ws.Range("T27").Value = 1
$ws.Range("T28").Value = 1
$ws.Range("T29").Value = 1
$ws.Range("T30").Value = 1
$ws.Range("T31").Value = 1
$ws.Range("T38").Value = 0
$ws.Range("T39").Value = 0
$ws.Range("T40").Value = 1
$ws.Range("T41").Value = 0
$ws.Range("T42").Value = 0
$ws.Range("T43").Value = 1
$ws.Range("T44").Value = 1
$ws.Range("T45").Value = 1
$ws.Range("T46").Value = 1
$ws.Range("T47").Value = 1
$ws.Range("T48").Value = 1
$ws.Range("T49").Value = 1
$ws.Range("T50").Value = 0
$ws.Range("T51").Value = 1
$ws.Range("T52").Value = 1
$ws.Range("T53").Value = 1
$ws.Range("T54").Value = 1
$ws.Range("T55").Value = 1
$ws.Range("T62").Value = 1
$ws.Range("T63").Value = 0
$ws.Range("T65").Value = 1
$ws.Range("T66").Value = 1
$ws.Range("T67").Value = 1
$ws.Range("T75").Value = 1
$ws.Range("T77").Value = 0
$ws.Range("T89").Value = 0
$ws.Range("T92").Value = 0
$ws.Range("T93").Value = 1
$ws.Range("T94").Value = 1
$ws.Range("T95").Value = 1
$ws.Range("T96").Value = 1
$ws.Range("T97").Value = 1
$ws.Range("T98").Value = 1
$ws.Range("T99").Value = 1
$ws.Range("T100").Value = 1
$ws.Range("T101").Value = 0
$ws.Range("T102").Value = 1
$ws.Range("T103").Value = 1
$ws.Range("T104").Value = 1
$ws.Range("T105").Value = 1
$ws.Range("T106").Value = 0
$ws.Range("T107").Value = 1
$ws.Range("T108").Value = 1
$ws.Range("T109").Value = 1
$ws.Range("T116").Value = 0
$ws.Range("T117").Value = 0
$ws.Range("T120").Value = 1
$ws.Range("T122").Value = 0
$ws.Range("T126").Value = 0
$ws.Range("T127").Value = 0
$ws.Range("T130").Value = 7
$ws.Range("T132").Value = 0
$ws.Range("T136").Value = 0
$ws.Range("T138").Value = 0
$ws.Range("T140").Value = 0
$ws.Range("T142").Value = 0
$ws.Range("T143").Value = 0
$ws.Range("T146").Value = 1
$ws.Range("T147").Value = 1
$ws.Range("T148").Value = 1
$ws.Range("T149").Value = 1
$ws.Range("T150").Value = 1
$ws.Range("T151").Value = 1
$ws.Range("T152").Value = 1
$ws.Range("T153").Value = 1
$ws.Range("T154").Value = 1
$ws.Range("T155").Value = 1
$ws.Range("T156").Value = 1
$ws.Range("T157").Value = 1
$ws.Range("T158").Value = 1
$ws.Range("T159").Value = 1
$ws.Range("T160").Value = 1
$ws.Range("T161").Value = 1
$ws.Range("T162").Value = 1
$ws.Range("T163").Value = 1
$ws.Range("T164").Value = 1
$ws.Range("T165").Value = 1
$ws.Range("T166").Value = 1
$ws.Range("T167").Value = 1
$ws.Range("T168").Value = 1
$ws.Range("T169").Value = 1
$ws.Range("T170").Value = 1
$ws.Range("T171").Value = 1
$ws.Range("T172").Value = 1
$ws.Range("T173").Value = 1
$ws.Range("T174").Value = 1
$ws.Range("T176").Value = 1
$ws.Range("T177").Value = 1
$ws.Range("T178").Value = 1
$ws.Range("T182").Value = 1
$ws.Range("T183").Value = 1
$ws.Range("T184").Value = 1
$ws.Range("T185").Value = 1
$ws.Range("T186").Value = 1
$ws.Range("T187").Value = 1
$ws.Range("T188").Value = 7
$ws.Range("T189").Value = 1
$ws.Range("T190").Value = 1
$ws.Range("T191").Value = 1
$ws.Range("T192").Value = 1
$ws.Range("T194").Value = 1
$ws.Range("T195").Value = 1
$ws.Range("T196").Value = 1
$ws.Range("T197").Value = 1
$ws.Range("T198").Value = 1
$ws.Range("T199").Value = 1
$ws.Range("T200").Value = 0
$ws.Range("T201").Value = 1
$ws.Range("T212").Value = 0
$ws.Range("T216").Value = 1
$ws.Range("T217").Value = 1
$ws.Range("T218").Value = 1
$ws.Range("T219").Value = 1
$ws.Range("T220").Value = 1
$ws.Range("T221").Value = 1
$ws.Range("T222").Value = 1
$ws.Range("T223").Value = 1
$ws.Range("T224").Value = 1
$ws.Range("T225").Value = 1
$ws.Range("T226").Value = 1
$ws.Range("T227").Value = 1
$ws.Range("T228").Value = 1
$ws.Range("T229").Value = 1
$ws.Range("T230").Value = 1
$ws.Range("T231").Value = 1
$ws.Range("T232").Value = 1
$ws.Range("T233").Value = 1
$ws.Range("T234").Value = 1
$ws.Range("T235").Value = 1
$ws.Range("T236").Value = 1
$ws.Range("T237").Value = 1
$ws.Range("T238").Value = 1
$ws.Range("T239").Value = 1
$ws.Range("T240").Value = 1
$ws.Range("T241").Value = 1
$ws.Range("T242").Value = 1
$ws.Range("T243").Value = 1
$ws.Range("T244").Value = 1
$ws.Range("T245").Value = 6
$ws.Range("T246").Value = 1
$ws.Range("T247").Value = 1
$ws.Range("T248").Value = 1
$ws.Range("T249").Value = 1
$ws.Range("T250").Value = 1
$ws.Range("T251").Value = 1
$ws.Range("T252").Value = 1
$ws.Range("T253").Value = 1
$ws.Range("T261").Value = 7
$ws.Range("T262").Value = 1
$ws.Range("T263").Value = 1
$ws.Range("T264").Value = 1
$ws.Range("T265").Value = 7
$ws.Range("T272").Value = 1
$ws.Range("T273").Value = 1
$ws.Range("T274").Value = 1
$ws.Range("T275").Value = 1
$ws.Range("T276").Value = 1
$ws.Range("T277").Value = 1
$ws.Range("T278").Value = 1
$ws.Range("T279").Value = 1
$ws.Range("T281").Value = 1
$ws.Range("T282").Value = 1
$ws.Range("T283").Value = 1
$ws.Range("T284").Value = 1
$ws.Range("T287").Value = 1
$ws.Range("T288").Value = 1
$ws.Range("T289").Value = 1
$ws.Range("T290").Value = 1
$ws.Range("T291").Value = 1
$ws.Range("T292").Value = 1
$ws.Range("T293").Value = 1
$ws.Range("T294").Value = 1
$ws.Range("T295").Value = 1
$ws.Range("T296").Value = 1
$ws.Range("T298").Value = 1
$ws.Range("T299").Value = 1
$ws.Range("T300").Value = 1
$ws.Range("T301").Value = 1
$ws.Range("T302").Value = 1
$ws.Range("T303").Value = 1
$ws.Range("T304").Value = 1
$ws.Range("T305").Value = 1
$ws.Range("T306").Value = 1
$ws.Range("T307").Value = 1
$ws.Range("T308").Value = 1
$ws.Range("T309").Value = 1
$ws.Range("T310").Value = 1
$ws.Range("T311").Value = 1
$ws.Range("T312").Value = 1
$ws.Range("T313").Value = 1
$ws.Range("T314").Value = 1
$ws.Range("T315").Value = 1
$ws.Range("T316").Value = 1
$ws.Range("T317").Value = 1
$ws.Range("T318").Value = 1
$ws.Range("T319").Value = 1
$ws.Range("T320").Value = 1
$ws.Range("T321").Value = 1
$ws.Range("T322").Value = 1
$ws.Range("T323").Value = 1
$ws.Range("T324").Value = 1
$ws.Range("T325").Value = 1
$ws.Range("T326").Value = 1
$ws.Range("T327").Value = 1
$ws.Range("T328").Value = 1
$ws.Range("T329").Value = 1
$ws.Range("T330").Value = 1
$ws.Range("T331").Value = 1
$ws.Range("T332").Value = 1
$ws.Range("T333").Value = 1
$ws.Range("T334").Value = 1
$ws.Range("T335").Value = 1
$ws.Range("T336").Value = 1
$ws.Range("T337").Value = 1
$ws.Range("T338").Value = 1
$ws.Range("T339").Value = 1
$ws.Range("T340").Value = 1
$ws.Range("T341").Value = 1
$ws.Range("T342").Value = 1
$ws.Range("T343").Value = 1
$ws.Range("T344").Value = 1
$ws.Range("T345").Value = 1
$ws.Range("T346").Value = 1
$ws.Range("T347").Value = 1
$ws.Range("T348").Value = 1
$ws.Range("T349").Value = 1
$ws.Range("T350").Value = 1
$ws.Range("T351").Value = 1
$ws.Range("T352").Value = 1
$ws.Range("T353").Value = 1
$ws.Range("T354").Value = 1
$ws.Range("T355").Value = 1
$ws.Range("T356").Value = 1
$ws.Range("T357").Value = 1
$ws.Range("T358").Value = 1
$ws.Range("T359").Value = 1
$ws.Range("T360").Value = 1
$ws.Range("T361").Value = 1
$ws.Range("T362").Value = 1
$ws.Range("T368").Value = 1
$ws.Range("T371").Value = 1
$ws.Range("T372").Value = 0
$ws.Range("T373").Value = 0
$ws.Range("T375").Value = 0
$ws.Range("T376").Value = 1
$ws.Range("T380").Value = 0
$ws.Range("T381").Value = 1
$ws.Range("T382").Value = 1
$ws.Range("T385").Value = 1
$ws.Range("T386").Value = 1
$ws.Range("T388").Value = 1
$ws.Range("T390").Value = 1
$ws.Range("T391").Value = 0
$ws.Range("T392").Value = 1
$ws.Range("T394").Value = 0
$ws.Range("T400").Value = 1
$ws.Range("T402").Value = 1
$ws.Range("T406").Value = 7
$ws.Range("T409").Value = 7
$ws.Range("T410").Value = 7
$ws.Range("T416").Value = 0
$ws.Range("T418").Value = 1
$ws.Range("T421").Value = 1
$ws.Range("T425").Value = 0
$ws.Range("T427").Value = 1
$ws.Range("T428").Value = 0
$ws.Range("T430").Value = 0
$ws.Range("T432").Value = 0
$ws.Range("T434").Value = 1
$ws.Range("T435").Value = 1
$ws.Range("T436").Value = 1
$ws.Range("T437").Value = 1
$ws.Range("T438").Value = 1
$ws.Range("T439").Value = 6
$ws.Range("T440").Value = 1
$ws.Range("T441").Value = 1
$ws.Range("T442").Value = 1
$ws.Range("T443").Value = 1
$ws.Range("T444").Value = 1
$ws.Range("T445").Value = 1
$ws.Range("T449").Value = 1
$ws.Range("T450").Value = 1
